# Add two new localization rows (the "red herring" joke lines) to Sheet1,
# following the existing key/value row pattern used throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "red_herring1"
$ws.Range("B14").Value = "If you’re at all familiar with logical fallacies, you’d realize what we’re looking at here is clearly a red herring – something to divert attention away from the real issue at hand."

$ws.Range("A15").Value = "red_herring2"
$ws.Range("B15").Value = "What is the real issue at hand?  I’m not sure.  I’m too distracted by the red fish flopping around."

# Leave the selection where the author's cursor ended up after adding the
# two new rows (two rows below the old A16 "next empty row" position).
$ws.Range("A18").Select() | Out-Null
